$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-text storage for numeric-looking Price cells so Excel keeps them as literal strings
# (matches the source file's inlineStr cell type) instead of silently converting to numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = '27.361.15'
$ws.Range("E2").Value = '  +1.37%  '
$ws.Range("D3").Value = '1.827.03'
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '314.23'
$ws.Range("E5").Value = '  +1.20%  '
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").Value = '0.4479'
$ws.Range("E7").Value = '  +5.51%  '
$ws.Range("D8").Value = '0.3773'
$ws.Range("E8").Value = '  +3.25%  '
$ws.Range("D9").Value = '0.07533'
$ws.Range("E9").Value = '  +4.32%  '
$ws.Range("D10").Value = '0.8935'
$ws.Range("E10").Value = '  +6.19%  '
$ws.Range("D11").Value = '21.09'
$ws.Range("E11").Value = '  +2.69%  '
$ws.Range("D12").Value = '1.815.68'
$ws.Range("E12").Value = '  -0.89%  '
$ws.Range("D13").Value = '6.747'
$ws.Range("E13").Value = '  +1.57%  '
$ws.Range("D14").Value = '94.60'
$ws.Range("E14").Value = '  +5.71%  '
$ws.Range("D15").Value = '5.413'
$ws.Range("E15").Value = '  +2.70%  '
$ws.Range("D16").Value = '0.07128'
$ws.Range("E16").Value = '  +0.87%  '
$ws.Range("E17").Value = '  -0.03%  '
$ws.Range("D18").Value = '0.000008837'
$ws.Range("E18").Value = '  +1.09%  '
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").Value = '15.23'
$ws.Range("E20").Value = '  +2.63%  '
$ws.Range("D21").Value = '27.378.25'
$ws.Range("E21").Value = '  +1.26%  '
$ws.Range("D22").Value = '5.288'
$ws.Range("E22").Value = '  +3.33%  '
$ws.Range("E23").Value = '  +1.65%  '
$ws.Range("D24").Value = '2.004'
$ws.Range("E24").Value = '  +1.59%  '
$ws.Range("D25").Value = '2.490'
$ws.Range("E25").Value = '  +11.85%  '
$ws.Range("E26").Value = '  +0.56%  '
$ws.Range("D27").Value = '18.66'
$ws.Range("E27").Value = '  +2.53%  '
$ws.Range("D28").Value = '5.374'
$ws.Range("E28").Value = '  +3.27%  '
$ws.Range("D29").Value = '118.06'
$ws.Range("E29").Value = '  +1.31%  '
$ws.Range("D30").Value = '0.08846'
$ws.Range("E30").Value = '  +1.58%  '
$ws.Range("D31").Value = '0.7822'
$ws.Range("E31").Value = '  +6.59%  '
$ws.Range("D32").Value = '1.206'
$ws.Range("E32").Value = '  +2.73%  '
$ws.Range("D33").Value = '4.565'
$ws.Range("E33").Value = '  +3.58%  '
$ws.Range("D34").Value = '2.890'
$ws.Range("E34").Value = '  -0.36%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D37").Value = '0.01993'
$ws.Range("E37").Value = '  +3.04%  '
$ws.Range("D38").Value = '0.05338'
$ws.Range("E38").Value = '  +2.55%  '
$ws.Range("E39").Value = '  +2.58%  '
$ws.Range("D40").Value = '0.5332'
$ws.Range("E40").Value = '  +4.46%  '
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").Value = '0.1733'
$ws.Range("E41").Value = '  +2.99%  '
$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").Value = '2.876'
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("D43").Value = '2.296'
$ws.Range("E43").Value = '  +17.63%  '
$ws.Range("E44").Value = '  +3.67%  '
$ws.Range("D45").Value = '0.5169'
$ws.Range("E45").Value = '  +9.38%  '
$ws.Range("D46").Value = '10.76'
$ws.Range("E46").Value = '  +2.17%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '1.712'
$ws.Range("E47").Value = '  +4.04%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '106.40'
$ws.Range("E48").Value = '  +0.64%  '
$ws.Range("E49").Value = '  +0.07%  '
$ws.Range("E50").Value = '  +0.92%  '
$ws.Range("D51").Value = '64.52'
$ws.Range("E51").Value = '  +3.42%  '
